$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("omnidirectional")
$ws.Range("C2").Value = 16.4
$ws.Range("D2").Value = 17.3
$ws.Range("E2").Value = 15.7
$ws.Range("C3").Value = 17.6
$ws.Range("D3").Value = 18.7
$ws.Range("E3").Value = 16.6
$ws.Range("C4").Value = 20.2
$ws.Range("D4").Value = 22.1
$ws.Range("E4").Value = 18.7

$ws = $wb.Worksheets.Item("345 to 15.0")
$ws.Range("C2").Value = 11.2
$ws.Range("D2").Value = 11.9
$ws.Range("E2").Value = 10.6
$ws.Range("C3").Value = 12.2
$ws.Range("D3").Value = 13.2
$ws.Range("E3").Value = 11.4
$ws.Range("C4").Value = 14.3
$ws.Range("D4").Value = 15.8
$ws.Range("E4").Value = 13.1

$ws = $wb.Worksheets.Item("15.0 to 45.0")
$ws.Range("C2").Value = 8.390000000000001
$ws.Range("D2").Value = 8.9
$ws.Range("E2").Value = 7.94
$ws.Range("C3").Value = 9.119999999999999
$ws.Range("D3").Value = 9.82
$ws.Range("E3").Value = 8.529999999999999
$ws.Range("C4").Value = 10.7
$ws.Range("D4").Value = 11.9
$ws.Range("E4").Value = 9.82

$ws = $wb.Worksheets.Item("45.0 to 75.0")
$ws.Range("C2").Value = 7.82
$ws.Range("D2").Value = 8.26
$ws.Range("E2").Value = 7.45
$ws.Range("C3").Value = 8.449999999999999
$ws.Range("D3").Value = 9.050000000000001
$ws.Range("E3").Value = 7.96
$ws.Range("C4").Value = 9.81
$ws.Range("D4").Value = 10.8
$ws.Range("E4").Value = 9.02

$ws = $wb.Worksheets.Item("75.0 to 105.0")
$ws.Range("C2").Value = 8.02
$ws.Range("D2").Value = 8.56
$ws.Range("E2").Value = 7.55
$ws.Range("C3").Value = 8.75
$ws.Range("D3").Value = 9.48
$ws.Range("E3").Value = 8.15
$ws.Range("C4").Value = 10.4
$ws.Range("D4").Value = 11.5
$ws.Range("E4").Value = 9.460000000000001

$ws = $wb.Worksheets.Item("105.0 to 135.0")
$ws.Range("C2").Value = 7.62
$ws.Range("D2").Value = 8.06
$ws.Range("E2").Value = 7.21
$ws.Range("C3").Value = 8.26
$ws.Range("D3").Value = 8.869999999999999
$ws.Range("E3").Value = 7.75
$ws.Range("C4").Value = 9.68
$ws.Range("D4").Value = 10.7
$ws.Range("E4").Value = 8.880000000000001

$ws = $wb.Worksheets.Item("135.0 to 165.0")
$ws.Range("C2").Value = 7.31
$ws.Range("D2").Value = 7.78
$ws.Range("E2").Value = 6.93
$ws.Range("C3").Value = 7.92
$ws.Range("D3").Value = 8.56
$ws.Range("E3").Value = 7.41
$ws.Range("C4").Value = 9.289999999999999
$ws.Range("D4").Value = 10.3
$ws.Range("E4").Value = 8.5

$ws = $wb.Worksheets.Item("165.0 to 195.0")
$ws.Range("C2").Value = 9.57
$ws.Range("D2").Value = 10.2
$ws.Range("E2").Value = 9.050000000000001
$ws.Range("C3").Value = 10.4
$ws.Range("D3").Value = 11.3
$ws.Range("E3").Value = 9.74
$ws.Range("C4").Value = 12.4
$ws.Range("D4").Value = 13.7
$ws.Range("E4").Value = 11.3

$ws = $wb.Worksheets.Item("195.0 to 225.0")
$ws.Range("C2").Value = 12.7
$ws.Range("D2").Value = 13.5
$ws.Range("E2").Value = 11.9
$ws.Range("C3").Value = 13.8
$ws.Range("D3").Value = 14.9
$ws.Range("E3").Value = 12.8
$ws.Range("C4").Value = 16.3
$ws.Range("D4").Value = 18
$ws.Range("E4").Value = 14.8

$ws = $wb.Worksheets.Item("225.0 to 255.0")
$ws.Range("C2").Value = 15.1
$ws.Range("D2").Value = 16.1
$ws.Range("E2").Value = 14.2
$ws.Range("C3").Value = 16.4
$ws.Range("D3").Value = 17.7
$ws.Range("E3").Value = 15.4
$ws.Range("C4").Value = 19.4
$ws.Range("D4").Value = 21.6
$ws.Range("E4").Value = 17.8

$ws = $wb.Worksheets.Item("255.0 to 285.0")
$ws.Range("C2").Value = 14.9
$ws.Range("D2").Value = 15.8
$ws.Range("E2").Value = 14.1
$ws.Range("C3").Value = 16.3
$ws.Range("D3").Value = 17.5
$ws.Range("E3").Value = 15.2
$ws.Range("C4").Value = 19.3
$ws.Range("D4").Value = 21.4
$ws.Range("E4").Value = 17.6

$ws = $wb.Worksheets.Item("285.0 to 315.0")
$ws.Range("C2").Value = 14.3
$ws.Range("D2").Value = 15.2
$ws.Range("E2").Value = 13.5
$ws.Range("C3").Value = 15.5
$ws.Range("D3").Value = 16.7
$ws.Range("E3").Value = 14.5
$ws.Range("C4").Value = 18.1
$ws.Range("D4").Value = 20.1
$ws.Range("E4").Value = 16.5

$ws = $wb.Worksheets.Item("315.0 to 345.0")
$ws.Range("C2").Value = 14.2
$ws.Range("D2").Value = 14.9
$ws.Range("E2").Value = 13.5
$ws.Range("C3").Value = 15.2
$ws.Range("D3").Value = 16.3
$ws.Range("E3").Value = 14.4
$ws.Range("C4").Value = 17.5
$ws.Range("D4").Value = 19.3
$ws.Range("E4").Value = 16.2
